# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Sun Apr 28 11:37:22 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + '63.672.00'
$ws.Range("E2").Value = '  +1.19%  '
# Row 3
$ws.Range("D3").Value = "'" + '3.283.66'
$ws.Range("E3").Value = '  +4.97%  '
# Row 4
$ws.Range("E4").Value = '  +0.02%  '
# Row 5
$ws.Range("D5").Value = "'" + '603.66'
$ws.Range("E5").Value = '  +2.61%  '
# Row 6
$ws.Range("D6").Value = "'" + '141.61'
$ws.Range("E6").Value = '  +4.26%  '
# Row 7
$ws.Range("E7").Value = '  +0.03%  '
# Row 8
$ws.Range("D8").Value = "'" + '3.280.22'
# Row 9
$ws.Range("E9").Value = '  +1.02%  '
# Row 10
$ws.Range("E10").Value = '  +3.33%  '
# Row 11
$ws.Range("E11").Value = '  +3.63%  '
# Row 12
$ws.Range("E12").Value = '  +3.54%  '
# Row 13
$ws.Range("E13").Value = '  +1.73%  '
# Row 14
$ws.Range("E14").Value = '  +1.72%  '
# Row 15
$ws.Range("D15").Value = "'" + '3.823.77'
$ws.Range("E15").Value = '  +5.18%  '
# Row 16
$ws.Range("E16").Value = '  +1.21%  '
# Row 17
$ws.Range("D17").Value = "'" + '3.282.62'
$ws.Range("E17").Value = '  +5.20%  '
# Row 18
$ws.Range("D18").Value = "'" + '63.716.75'
$ws.Range("E18").Value = '  +1.14%  '
# Row 19
$ws.Range("E19").Value = '  +3.29%  '
# Row 20
$ws.Range("D20").Value = "'" + '479.31'
$ws.Range("E20").Value = '  +2.13%  '
# Row 21
$ws.Range("E21").Value = '  +0.01%  '
# Row 22
$ws.Range("E22").Value = '  +4.83%  '
# Row 23
$ws.Range("E23").Value = '  +4.59%  '
# Row 24
$ws.Range("D24").Value = "'" + '13.49'
$ws.Range("E24").Value = '  +4.63%  '
# Row 25
$ws.Range("D25").Value = "'" + '84.20'
$ws.Range("E25").Value = '  -1.14%  '
# Row 26
$ws.Range("D26").Value = "'" + '0.999'
$ws.Range("E26").Value = '  -0.17%  '
# Row 27
$ws.Range("E27").Value = '  +2.61%  '
# Row 28
$ws.Range("D28").Value = "'" + '7.30'
$ws.Range("E28").Value = '  +7.31%  '
# Row 29
$ws.Range("E29").Value = '  +0.00%  '
# Row 30
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = "'" + '8.10'
$ws.Range("E30").Value = '  +3.32%  '
# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'" + '2.16'
$ws.Range("E31").Value = '  +3.95%  '
# Row 32
$ws.Range("D32").Value = "'" + '28.79'
$ws.Range("E32").Value = '  +8.46%  '
# Row 33
$ws.Range("E33").Value = '  -3.16%  '
# Row 34
$ws.Range("E34").Value = '  +0.20%  '
# Row 35
$ws.Range("D35").Value = "'" + '1.10'
$ws.Range("E35").Value = '  +3.16%  '
# Row 36
$ws.Range("E36").Value = '  +4.07%  '
# Row 37
$ws.Range("D37").Value = "'" + '52.99'
$ws.Range("E37").Value = '  +2.00%  '
# Row 38
$ws.Range("D38").Value = "'" + '0.0₃0739'
$ws.Range("E38").Value = '  +9.51%  '
# Row 39
$ws.Range("D39").Value = "'" + '0.0397'
$ws.Range("E39").Value = '  +3.40%  '
# Row 40
$ws.Range("D40").Value = "'" + '425.23'
$ws.Range("E40").Value = '  +2.19%  '
# Row 41
$ws.Range("D41").Value = "'" + '3.054.64'
$ws.Range("E41").Value = '  +5.05%  '
# Row 42
$ws.Range("E42").Value = '  +2.26%  '
# Row 43
$ws.Range("E43").Value = '  +2.27%  '
# Row 44
$ws.Range("E44").Value = '  +1.30%  '
# Row 45
$ws.Range("E45").Value = '  +2.54%  '
# Row 46
$ws.Range("E46").Value = '  +4.40%  '
# Row 47
$ws.Range("D47").Value = "'" + '26.16'
$ws.Range("E47").Value = '  +3.71%  '
# Row 48
$ws.Range("E48").Value = '  +0.03%  '
# Row 50
$ws.Range("D50").Value = "'" + '124.56'
$ws.Range("E50").Value = '  +3.36%  '
# Row 51
$ws.Range("E51").Value = '  +2.16%  '
